$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 351.55554
$ws.Range("I9").Value = 72.8
$ws.Range("J9").Value = 700
$ws.Range("K9").Value = 72.8
$ws.Range("L9").Value = 700
$ws.Range("M9").Value = 96.2
$ws.Range("N9").Value = -1038
$ws.Range("H116").Value = 4431.579
$ws.Range("I116").Value = 3853.3333
$ws.Range("J116").Value = 4952
$ws.Range("K116").Value = 3853.3333
$ws.Range("L116").Value = 4952
$ws.Range("M116").Value = -411.3332999999998
$ws.Range("N116").Value = -11836
$ws.Range("H127").Value = 1651.1111
$ws.Range("J127").Value = 1545
$ws.Range("L127").Value = 4635
$ws.Range("N127").Value = -14555
$ws.Range("H129").Value = 3732314.5
$ws.Range("J129").Value = 924.873
$ws.Range("L129").Value = 2774.619
$ws.Range("N129").Value = -12774.619
$ws.Range("H131").Value = 2277.6128
$ws.Range("I131").Value = 2277.0588
$ws.Range("J131").Value = 2278.2856
$ws.Range("K131").Value = 6831.176399999999
$ws.Range("L131").Value = 6834.8568
$ws.Range("M131").Value = -1791.176399999999
$ws.Range("N131").Value = -16914.8568
$ws.Range("H135").Value = 487.5
$ws.Range("I135").Value = 358.30554
$ws.Range("J135").Value = 1650.25
$ws.Range("K135").Value = 3224.74986
$ws.Range("L135").Value = 14852.25
$ws.Range("M135").Value = -689.7498599999999
$ws.Range("N135").Value = -19922.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7137.11
$ws.Range("I32").Value = 7137.11
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 7137.11
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -6850.11
$ws.Range("N32").ClearContents()
$ws.Range("H122").Value = 3600.6956
$ws.Range("I122").Value = 2773.4666
$ws.Range("J122").Value = 5151.75
$ws.Range("K122").Value = 8320.399800000001
$ws.Range("L122").Value = 15455.25
$ws.Range("M122").Value = -5870.399800000001
$ws.Range("N122").Value = -20355.25
$ws.Range("H132").Value = 2150.1785
$ws.Range("I132").Value = 1795.625
$ws.Range("J132").Value = 4277.5
$ws.Range("K132").Value = 5386.875
$ws.Range("L132").Value = 12832.5
$ws.Range("M132").Value = -2856.875
$ws.Range("N132").Value = -17892.5
$ws.Range("H137").Value = 29360
$ws.Range("J137").Value = 29360
$ws.Range("L137").Value = 29360
$ws.Range("N137").Value = -39560

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 13530.875
$ws.Range("I75").Value = 5414.25
$ws.Range("J75").Value = 21647.5
$ws.Range("K75").Value = 5414.25
$ws.Range("L75").Value = 21647.5
$ws.Range("M75").Value = -4478.25
$ws.Range("N75").Value = -23519.5
$ws.Range("H78").Value = 13530.875
$ws.Range("I78").Value = 5414.25
$ws.Range("J78").Value = 21647.5
$ws.Range("K78").Value = 16242.75
$ws.Range("L78").Value = 64942.5
$ws.Range("M78").Value = -11562.75
$ws.Range("N78").Value = -74302.5
$ws.Range("H99").Value = 3934.8333
$ws.Range("I99").Value = 3468.6667
$ws.Range("K99").Value = 3468.6667
$ws.Range("M99").Value = -1970.6667
$ws.Range("H107").Value = 1768.36
$ws.Range("I107").Value = 1459.8
$ws.Range("J107").Value = 2231.2
$ws.Range("K107").Value = 1459.8
$ws.Range("L107").Value = 2231.2
$ws.Range("M107").Value = 460.2
$ws.Range("N107").Value = -6071.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 23818.6
$ws.Range("J51").Value = 23818.6
$ws.Range("L51").Value = 23818.6
$ws.Range("N51").Value = -25290.6
$ws.Range("H61").Value = 23818.6
$ws.Range("J61").Value = 23818.6
$ws.Range("L61").Value = 23818.6
$ws.Range("N61").Value = -24514.6
$ws.Range("H132").Value = 2498.4243
$ws.Range("I132").Value = 2124.074
$ws.Range("K132").Value = 6372.222
$ws.Range("M132").Value = -3842.222
$ws.Range("H141").Value = 25718.75
$ws.Range("J141").Value = 25718.75
$ws.Range("L141").Value = 25718.75
$ws.Range("N141").Value = -36078.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 1260
$ws.Range("I20").Value = 1260
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 3780
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -3553
$ws.Range("N20").ClearContents()
$ws.Range("H40").Value = 68.666664
$ws.Range("I40").Value = 68.666664
$ws.Range("K40").Value = 274.666656
$ws.Range("M40").Value = -205.666656
$ws.Range("H93").Value = 2779.8
$ws.Range("H94").Value = 3875.862
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 3942.8572
$ws.Range("K94").Value = 6000
$ws.Range("L94").Value = 11828.5716
$ws.Range("M94").Value = -5324
$ws.Range("N94").Value = -13180.5716
$ws.Range("H95").Value = 3854.2856
$ws.Range("I95").Value = 3000
$ws.Range("J95").Value = 3996.6667
$ws.Range("K95").Value = 9000
$ws.Range("L95").Value = 11990.0001
$ws.Range("M95").Value = -6941
$ws.Range("N95").Value = -16108.0001
$ws.Range("H96").Value = 4000
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H99").Value = 1937.5
$ws.Range("I99").Value = 875
$ws.Range("K99").Value = 2625
$ws.Range("M99").Value = -379
$ws.Range("H102").Value = 2884.7917
$ws.Range("I102").Value = 1777.5
$ws.Range("J102").Value = 2985.4546
$ws.Range("K102").Value = 5332.5
$ws.Range("L102").Value = 8956.363799999999
$ws.Range("M102").Value = -2898.5
$ws.Range("N102").Value = -13824.3638
$ws.Range("H108").Value = 2891.889
$ws.Range("I108").Value = 2027
$ws.Range("K108").Value = 6081
$ws.Range("M108").Value = -3201
$ws.Range("H110").Value = 3250
$ws.Range("H111").Value = 563.5
$ws.Range("I111").Value = 563.5
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 1690.5
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("M111").Value = 1376.5
$ws.Range("H131").Value = 1413.2456
$ws.Range("I131").Value = 2730.9
$ws.Range("J131").Value = 1132.8937
$ws.Range("K131").Value = 8192.700000000001
$ws.Range("L131").Value = 3398.6811
$ws.Range("M131").Value = -3152.700000000001
$ws.Range("N131").Value = -13478.6811

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 930
$ws.Range("I107").Value = 405.30768
$ws.Range("K107").Value = 405.30768
$ws.Range("M107").Value = 1514.69232
$ws.Range("H132").Value = 4216
$ws.Range("I132").Value = 4463.154
$ws.Range("J132").Value = 3757
$ws.Range("K132").Value = 13389.462
$ws.Range("L132").Value = 11271
$ws.Range("M132").Value = -10859.462
$ws.Range("N132").Value = -16331

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3307.7693
$ws.Range("I122").Value = 2720.4
$ws.Range("J122").Value = 3674.875
$ws.Range("K122").Value = 8161.200000000001
$ws.Range("L122").Value = 11024.625
$ws.Range("M122").Value = -5711.200000000001
$ws.Range("N122").Value = -15924.625
$ws.Range("H132").Value = 2504.0227
$ws.Range("I132").Value = 1763.4642
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 5290.392599999999
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -2760.392599999999
$ws.Range("N132").Value = -16460

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1876.8096
$ws.Range("I122").Value = 1478.5555
$ws.Range("J122").Value = 4266.3335
$ws.Range("K122").Value = 4435.666499999999
$ws.Range("L122").Value = 12799.0005
$ws.Range("M122").Value = -1985.666499999999
$ws.Range("N122").Value = -17699.0005
$ws.Range("H132").Value = 19201.96
$ws.Range("I132").Value = 3644.3684
$ws.Range("K132").Value = 10933.1052
$ws.Range("M132").Value = -8403.1052
